$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force columns D and E to text format so numeric-looking strings
# (e.g. "26.048.20", "40.58") are stored as text, matching the source data.
$rngD = $ws.Range("D2:D51")
$rngE = $ws.Range("E2:E51")
$rngD.NumberFormat = "@"
$rngE.NumberFormat = "@"

$ws.Range("D2").Value = "26.048.20"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "1.763.56"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "237.48"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.5263"
$ws.Range("E7").Value = "  +4.22%  "
$ws.Range("D8").Value = "0.2745"
$ws.Range("E8").Value = "  +3.26%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "40.58"
$ws.Range("E9").Value = "  -2.92%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.06214"
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.774.77"
$ws.Range("E11").Value = "  +1.80%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "16.07"
$ws.Range("E12").Value = "  +5.32%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.07038"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "0.6524"
$ws.Range("E14").Value = "  +9.49%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "4.519"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "78.46"
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "0.9999"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "0.9997"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "26.056.32"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "11.74"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.000006759"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "1.995.79"
$ws.Range("E22").Value = "  +1.58%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "4.099"
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "8.429"
$ws.Range("E24").Value = "  +3.36%  "
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").Value = "5.206"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "137.96"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "1.483"
$ws.Range("E27").Value = "  -2.35%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "1.846"
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "15.25"
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "102.99"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.08428"
$ws.Range("E31").Value = "  +3.95%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "3.731"
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "3.461"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.04449"
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.653"
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.003"
$ws.Range("E36").Value = "  +1.99%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "0.6134"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "2.757"
$ws.Range("E38").Value = "  +3.76%  "
$ws.Range("D39").Value = "1.994"
$ws.Range("E39").Value = "  +4.48%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.01586"
$ws.Range("E40").Value = "  +2.44%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "103.24"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "0.3912"
$ws.Range("E43").Value = "  +2.93%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "0.7559"
$ws.Range("E44").Value = "  +3.75%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "4.991"
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("D46").Value = "6.534"
$ws.Range("E46").Value = "  +11.00%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.05509"
$ws.Range("E47").Value = "  +3.24%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "0.1129"
$ws.Range("E48").Value = "  +1.49%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "30.27"
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "53.08"
$ws.Range("E50").Value = "  +1.10%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "0.3468"
$ws.Range("E51").Value = "  +0.68%  "

# Remove the temporary number-format styling so the resulting cells
# have no explicit style, same as the original workbook.
$rngD.ClearFormats()
$rngE.ClearFormats()
